# Daily attendance processing - 2025-11-28 16:31:11
# Normalize the "Recorded By" column (G): move any "System"/"system"
# entries to the end of the comma-separated list, preserving the
# relative order of the remaining entries (and of the System entries
# among themselves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = $text.Split(",")

    $others = @()
    $systems = @()

    foreach ($part in $parts) {
        $trimmedPart = $part.Trim()
        if ($trimmedPart -eq "System") {
            $systems += $trimmedPart
        } else {
            $others += $trimmedPart
        }
    }

    $newParts = $others + $systems
    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
